$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K shifts to F:M)
$ws.Range("D:E").Insert()

# Copy number/date formatting from column F (old D, now shifted) onto new D:E columns
$ws.Range("F7:F102").Copy()
$ws.Range("D7:E102").PasteSpecial(-4122)
$ws.Range("F38:F102").Copy()
$ws.Range("D38:E102").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("D80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Rows 37 and 79 are section headers with no data columns; undo the accidental
# formatting paste so these rows stay exactly as they were (label cell only)
$ws.Range("D37:E37").Clear()
$ws.Range("D79:E79").Clear()

# Populate the two new quarter columns (D = most recent quarter, E = prior quarter)
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43371
$ws.Range("D8").Value = 6944800
$ws.Range("E8").Value = 6710600
$ws.Range("D9").Value = 6527100
$ws.Range("E9").Value = 6308300
$ws.Range("D10").Value = 417700
$ws.Range("E10").Value = 402300
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 65800
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = 20300
$ws.Range("E15").Value = 18200
$ws.Range("D17").Value = 6922700
$ws.Range("E17").Value = 6560800
$ws.Range("D18").Value = 22100
$ws.Range("E18").Value = 149900
$ws.Range("D20").Value = -15300
$ws.Range("E20").Value = -6000
$ws.Range("D21").Value = 245000
$ws.Range("E21").Value = 291200
$ws.Range("D22").Value = 38800
$ws.Range("E22").Value = 35100
$ws.Range("D23").Value = -31900
$ws.Range("E23").Value = 108800
$ws.Range("D24").Value = 13300
$ws.Range("E24").Value = 21900
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = -45200
$ws.Range("E26").Value = 86900
$ws.Range("D27").Value = -45200
$ws.Range("E27").Value = 86900
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 15300
$ws.Range("E32").Value = 6000
$ws.Range("D33").Value = -45200
$ws.Range("E33").Value = 86900
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = -45200
$ws.Range("E35").Value = 86900
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43371
$ws.Range("D41").Value = 1503400
$ws.Range("E41").Value = 1377700
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 3160300
$ws.Range("E43").Value = 3277600
$ws.Range("D44").Value = 3897900
$ws.Range("E44").Value = 4442900
$ws.Range("D45").Value = 930400
$ws.Range("E45").Value = 935000
$ws.Range("D46").Value = 9491900
$ws.Range("E46").Value = 10033200
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 2214100
$ws.Range("E48").Value = 2277900
$ws.Range("D49").Value = 1428500
$ws.Range("E49").Value = 1457900
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 840500
$ws.Range("E52").Value = 957200
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 13975100
$ws.Range("E54").Value = 14726200
$ws.Range("D57").Value = 5543300
$ws.Range("E57").Value = 6236000
$ws.Range("D58").Value = 43200
$ws.Range("E58").Value = 55600
$ws.Range("D59").Value = 1917000
$ws.Range("E59").Value = 1862900
$ws.Range("D60").Value = 7503600
$ws.Range("E60").Value = 8154500
$ws.Range("D61").Value = 2906300
$ws.Range("E61").Value = 2869600
$ws.Range("D62").Value = 486900
$ws.Range("E62").Value = 532600
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 10896800
$ws.Range("E66").Value = 11556700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -2947700
$ws.Range("E72").Value = -2902500
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 3078300
$ws.Range("E76").Value = 3169500
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43371
$ws.Range("D81").Value = -45200
$ws.Range("E81").Value = 86900
$ws.Range("D83").Value = 238100
$ws.Range("E83").Value = 147300
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = -621000
$ws.Range("E89").Value = -1035800
$ws.Range("D91").Value = -228700
$ws.Range("E91").Value = -191100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = 768700
$ws.Range("E94").Value = 1230900
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -18100
$ws.Range("E100").Value = -62400
$ws.Range("D101").Value = -3900
$ws.Range("E101").Value = -9600
$ws.Range("D102").Value = 125600
$ws.Range("E102").Value = 123100


# Refresh financial figures that changed beyond the simple column shift
$ws.Range("I91").Value = -139200
$ws.Range("J91").Value = -124900

$ws.Range("G94").Value = -3664700
$ws.Range("H94").Value = 921100
$ws.Range("I94").Value = 1352700
